# Apply updated H:N profit-calculation values for the rows identified in the
# scheduled-runner refresh across all eight job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = $null

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = $null

$ws.Range("H111").Value = 2460
$ws.Range("I111").Value = 2460
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 7380
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -4313

$ws.Range("H113").Value = 3200
$ws.Range("I113").Value = 3200
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3200
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 54
$ws.Range("N113").Value = $null

$ws.Range("H125").Value = 4218.5713
$ws.Range("I125").Value = 2265.5
$ws.Range("J125").Value = 4999.8
$ws.Range("K125").Value = 20389.5
$ws.Range("L125").Value = 44998.2
$ws.Range("M125").Value = -17929.5
$ws.Range("N125").Value = -49918.2

$ws.Range("H132").Value = 1203.3478
$ws.Range("I132").Value = 1161.2858
$ws.Range("J132").Value = 1645
$ws.Range("K132").Value = 3483.8574
$ws.Range("L132").Value = 4935
$ws.Range("M132").Value = -953.8574000000003
$ws.Range("N132").Value = -9995

$ws.Range("H137").Value = 4435.2
$ws.Range("I137").Value = 2411.889
$ws.Range("J137").Value = 6090.636
$ws.Range("K137").Value = 7235.667
$ws.Range("L137").Value = 18271.908
$ws.Range("M137").Value = -4685.667
$ws.Range("N137").Value = -23371.908

$ws.Range("H138").Value = 2682.5715
$ws.Range("I138").Value = 1087.8214
$ws.Range("J138").Value = 4277.3213
$ws.Range("K138").Value = 3263.4642
$ws.Range("L138").Value = 12831.9639
$ws.Range("M138").Value = 1876.5358
$ws.Range("N138").Value = -23111.9639

$ws.Range("H141").Value = 2902.3333
$ws.Range("I141").Value = 2015.1875
$ws.Range("J141").Value = 9999.5
$ws.Range("K141").Value = 6045.5625
$ws.Range("L141").Value = 29998.5
$ws.Range("M141").Value = -865.5625
$ws.Range("N141").Value = -40358.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4546.375
$ws.Range("I32").Value = 3849.8333
$ws.Range("J32").Value = 14994.5
$ws.Range("K32").Value = 3849.8333
$ws.Range("L32").Value = 14994.5
$ws.Range("M32").Value = -3562.8333
$ws.Range("N32").Value = -15568.5

$ws.Range("H102").Value = 2626.0908
$ws.Range("I102").Value = 2138.7
$ws.Range("J102").Value = 7500
$ws.Range("K102").Value = 2138.7
$ws.Range("L102").Value = 7500
$ws.Range("M102").Value = -516.6999999999998
$ws.Range("N102").Value = -10744

$ws.Range("H122").Value = 1848.375
$ws.Range("I122").Value = 1848.375
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5545.125
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3095.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = $null

$ws.Range("H96").Value = 99999
$ws.Range("I96").Value = 99999
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 99999
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -97253
$ws.Range("N96").Value = $null

$ws.Range("H134").Value = 3582.8333
$ws.Range("I134").Value = 3439.8235
$ws.Range("J134").Value = 6014
$ws.Range("K134").Value = 10319.4705
$ws.Range("L134").Value = 18042
$ws.Range("M134").Value = -7784.470499999999
$ws.Range("N134").Value = -23112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3248.5
$ws.Range("I31").Value = 3248.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 3248.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2953.5

$ws.Range("H34").Value = 3248.5
$ws.Range("I34").Value = 3248.5
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 3248.5
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -3046.5

$ws.Range("H97").Value = 40000
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 40000
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 40000
$ws.Range("N97").Value = -41982

$ws.Range("H132").Value = 3427.3684
$ws.Range("I132").Value = 3187.2666
$ws.Range("J132").Value = 4327.75
$ws.Range("K132").Value = 9561.7998
$ws.Range("L132").Value = 12983.25
$ws.Range("M132").Value = -7031.799800000001
$ws.Range("N132").Value = -18043.25

$ws.Range("H134").Value = 5440.727
$ws.Range("I134").Value = 5440.727
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 16322.181
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -13787.181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 753.6
$ws.Range("I5").Value = 692
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 2076
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -1964
$ws.Range("N5").Value = -3224

$ws.Range("H75").Value = 500
$ws.Range("I75").Value = 500
$ws.Range("J75").Value = 500
$ws.Range("K75").Value = 1500
$ws.Range("L75").Value = 1500
$ws.Range("M75").Value = -502
$ws.Range("N75").Value = -3496

$ws.Range("H78").Value = 500
$ws.Range("I78").Value = 500
$ws.Range("J78").Value = 500
$ws.Range("K78").Value = 4500
$ws.Range("L78").Value = 4500
$ws.Range("M78").Value = 492
$ws.Range("N78").Value = -14484

$ws.Range("H81").Value = 706.5
$ws.Range("I81").Value = 706.5
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2119.5
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -996.5

$ws.Range("H84").Value = 706.5
$ws.Range("I84").Value = 706.5
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 6358.5
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -742.5

$ws.Range("H122").Value = 1226.75
$ws.Range("I122").Value = 703.5
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 6331.5
$ws.Range("L122").Value = 15750
$ws.Range("M122").Value = -3881.5
$ws.Range("N122").Value = -20650

$ws.Range("H135").Value = 753.6
$ws.Range("I135").Value = 692
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 6228
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -3693
$ws.Range("N135").Value = -14070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = $null

$ws.Range("H132").Value = 2666.6667
$ws.Range("I132").Value = 2666.6667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8000.000100000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5470.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2299.5557
$ws.Range("I7").Value = 2337
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 2337
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -2225
$ws.Range("N7").Value = -2224

$ws.Range("H68").Value = 3250
$ws.Range("I68").Value = 3250
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 3250
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2501

$ws.Range("H71").Value = 3250
$ws.Range("I71").Value = 3250
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 16250
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -12506

$ws.Range("H82").Value = 1557.7142
$ws.Range("I82").Value = 1461
$ws.Range("J82").Value = 1799.5
$ws.Range("K82").Value = 1461
$ws.Range("L82").Value = 1799.5
$ws.Range("M82").Value = -1100
$ws.Range("N82").Value = -2521.5

$ws.Range("H85").Value = 1557.7142
$ws.Range("I85").Value = 1461
$ws.Range("J85").Value = 1799.5
$ws.Range("K85").Value = 1461
$ws.Range("L85").Value = 1799.5
$ws.Range("M85").Value = -213
$ws.Range("N85").Value = -4295.5

$ws.Range("H126").Value = 2299.5557
$ws.Range("I126").Value = 2337
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 7011
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -4541
$ws.Range("N126").Value = -10940

$ws.Range("H136").Value = 2788.4375
$ws.Range("I136").Value = 2788.4375
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8365.3125
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5815.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5104
$ws.Range("I81").Value = 1129.5
$ws.Range("J81").Value = 21002
$ws.Range("K81").Value = 2259
$ws.Range("L81").Value = 42004
$ws.Range("M81").Value = -1198
$ws.Range("N81").Value = -44126

$ws.Range("H84").Value = 5104
$ws.Range("I84").Value = 1129.5
$ws.Range("J84").Value = 21002
$ws.Range("K84").Value = 11295
$ws.Range("L84").Value = 210020
$ws.Range("M84").Value = -5991
$ws.Range("N84").Value = -220628

$ws.Range("H100").Value = 891.6
$ws.Range("I100").Value = 617.6
$ws.Range("J100").Value = 1165.6
$ws.Range("K100").Value = 1235.2
$ws.Range("L100").Value = 2331.2
$ws.Range("M100").Value = -694.2
$ws.Range("N100").Value = -3413.2

$ws.Range("H122").Value = 2160
$ws.Range("I122").Value = 2160
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6480
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4030
$ws.Range("N122").Value = $null

$ws.Range("H126").Value = 2609.7144
$ws.Range("I126").Value = 2628.8333
$ws.Range("J126").Value = 2495
$ws.Range("K126").Value = 7886.499899999999
$ws.Range("L126").Value = 7485
$ws.Range("M126").Value = -5416.499899999999
$ws.Range("N126").Value = -12425

$ws.Range("H132").Value = 1379.619
$ws.Range("I132").Value = 1261.7059
$ws.Range("J132").Value = 1880.75
$ws.Range("K132").Value = 3785.1177
$ws.Range("L132").Value = 5642.25
$ws.Range("M132").Value = -1255.1177
$ws.Range("N132").Value = -10702.25

$ws.Range("H136").Value = 936.5769
$ws.Range("I136").Value = 839.4167
$ws.Range("J136").Value = 2102.5
$ws.Range("K136").Value = 2518.2501
$ws.Range("L136").Value = 6307.5
$ws.Range("M136").Value = 31.7498999999998
$ws.Range("N136").Value = -11407.5
